$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.597.24"
$ws.Range("E2").Value = "  -2.21%  "

$ws.Range("D3").Value = "'1.791.02"
$ws.Range("E3").Value = "  -2.12%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'231.50"
$ws.Range("E5").Value = "  -1.79%  "

$ws.Range("D6").Value = "'0.5880"
$ws.Range("E6").Value = "  -2.56%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").Value = "'0.2760"
$ws.Range("E8").Value = "  -1.28%  "

$ws.Range("D9").Value = "'0.06731"
$ws.Range("E9").Value = "  -4.64%  "

$ws.Range("D10").Value = "'23.12"
$ws.Range("E10").Value = "  -1.99%  "

$ws.Range("D11").Value = "'0.07527"
$ws.Range("E11").Value = "  -1.61%  "

$ws.Range("D12").Value = "'1.788.25"
$ws.Range("E12").Value = "  -2.40%  "

$ws.Range("D13").Value = "'4.783"
$ws.Range("E13").Value = "  -0.18%  "

$ws.Range("D14").Value = "'0.6120"
$ws.Range("E14").Value = "  -2.75%  "

$ws.Range("D15").Value = "'2.033.70"

$ws.Range("D16").Value = "'75.23"
$ws.Range("E16").Value = "  -4.87%  "

$ws.Range("D17").Value = "'0.000008882"
$ws.Range("E17").Value = "  -9.94%  "

$ws.Range("D18").Value = "'28.582.47"
$ws.Range("E18").Value = "  -2.23%  "

$ws.Range("D19").Value = "'5.402"
$ws.Range("E19").Value = "  -7.54%  "

$ws.Range("D21").Value = "'208.99"
$ws.Range("E21").Value = "  -6.81%  "

$ws.Range("D22").Value = "'11.45"
$ws.Range("E22").Value = "  -2.22%  "

$ws.Range("D23").Value = "'6.815"
$ws.Range("E23").Value = "  -2.76%  "

$ws.Range("D24").Value = "'1.005"
$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("D25").Value = "'152.82"
$ws.Range("E25").Value = "  -2.40%  "

$ws.Range("D26").Value = "'8.131"
$ws.Range("E26").Value = "  +1.87%  "

$ws.Range("E27").Value = "  -3.23%  "

$ws.Range("D28").Value = "'16.37"
$ws.Range("E28").Value = "  -1.51%  "

$ws.Range("D29").Value = "'1.415"
$ws.Range("E29").Value = "  -4.14%  "

$ws.Range("D30").Value = "'0.06201"
$ws.Range("E30").Value = "  -6.54%  "

$ws.Range("D31").Value = "'1.420"
$ws.Range("E31").Value = "  -1.94%  "

$ws.Range("D32").Value = "'3.805"
$ws.Range("E32").Value = "  +0.21%  "

$ws.Range("D33").Value = "'3.776"
$ws.Range("E33").Value = "  -1.85%  "

$ws.Range("D34").Value = "'1.734"
$ws.Range("E34").Value = "  +0.75%  "

$ws.Range("D35").Value = "'1.044"
$ws.Range("E35").Value = "  -5.74%  "

$ws.Range("D36").Value = "'0.6383"
$ws.Range("E36").Value = "  -1.36%  "

$ws.Range("D37").Value = "'2.503"
$ws.Range("E37").Value = "  -1.64%  "

$ws.Range("D38").Value = "'2.712"
$ws.Range("E38").Value = "  -0.96%  "

$ws.Range("D39").Value = "'6.395"
$ws.Range("E39").Value = "  -2.47%  "

$ws.Range("D40").Value = "'0.01694"
$ws.Range("E40").Value = "  -3.04%  "

$ws.Range("D41").Value = "'1.138.85"
$ws.Range("E41").Value = "  -6.23%  "

$ws.Range("D42").Value = "'0.8753"
$ws.Range("E42").Value = "  -2.53%  "

$ws.Range("D43").Value = "'1.004"
$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("D44").Value = "'99.99"
$ws.Range("E44").Value = "  -0.36%  "

$ws.Range("D45").Value = "'1.943.75"
$ws.Range("E45").Value = "  -2.53%  "

$ws.Range("D46").Value = "'59.79"
$ws.Range("E46").Value = "  -4.74%  "

$ws.Range("E47").Value = "  -4.89%  "

$ws.Range("D48").Value = "'1.582"
$ws.Range("E48").Value = "  -0.14%  "

$ws.Range("D49").Value = "'8.346"
$ws.Range("E49").Value = "  -2.59%  "

$ws.Range("D50").Value = "'0.05468"
$ws.Range("E50").Value = "  -0.52%  "

$ws.Range("D51").Value = "'0.4487"
$ws.Range("E51").Value = "  -1.36%  "
